$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$wsFE = $wb.Worksheets.Item("FE")
$wsTS = $wb.Worksheets.Item("TS")

# ---------------------------------------------------------------------
# Sheet1 ("Apuestas" summary table, Table4): add an "Estado" column and
# five new rows (FE09/TS05, FE10, FE11, FE12).
# ---------------------------------------------------------------------
$loApuestas = $ws1.ListObjects.Item("Table4")
$loApuestas.ListColumns.Add()
for ($i = 0; $i -lt 5; $i++) {
    $loApuestas.ListRows.Add()
}

# Header for the new column
$ws1.Range("E8").Value = "Estado"

# Remove the bold-ish style previously applied to D9 / D11
$ws1.Range("D9").ClearFormats()
$ws1.Range("D11").ClearFormats()
$ws1.Range("D9").Value = "ZZ01"
$ws1.Range("D11").Value = "ZZ01"

# New rows
$ws1.Range("A14").Value = "FE09"
$ws1.Range("F14").Value = "Parece que terminó, no 100% seguro"

$ws1.Range("A15").Value = "FE05"
$ws1.Range("B15").Value = "TS05"
$ws1.Range("F15").Value = "OK"

$ws1.Range("A16").Value = "FE10"
$ws1.Range("E16").Value = "corriendo"

$ws1.Range("A17").Value = "FE11"
$ws1.Range("E17").Value = "corriendo"

$ws1.Range("A18").Value = "FE12"
$ws1.Range("E18").Value = "corriendo"

# ---------------------------------------------------------------------
# FE sheet (Table1): row 11 FE08 -> FE09, plus three new rows
# (FE10/FE11/FE12).
# ---------------------------------------------------------------------
$loFE = $wsFE.ListObjects.Item("Table1")
$wsFE.Range("A11").Value = "FE09"

$loFE.ListRows.Add()
$loFE.ListRows.Add()
$loFE.ListRows.Add()

$wsFE.Range("A12").Value = "FE10"
$wsFE.Range("B12").Value = 300
$wsFE.Range("C12").Value = 11
$wsFE.Range("D12").Value = 400
$wsFE.Range("E12").Value = 40
$wsFE.Range("F12").Value = 10881
$wsFE.Range("G12").Value = 0.25
$wsFE.Range("H12").Value = $true
$wsFE.Range("I12").Value = $true
$wsFE.Range("J12").Value = $true
$wsFE.Range("K12").Value = "Todo TRUE"
$wsFE.Range("L12").Value = "Todo TRUE"

$wsFE.Range("A13").Value = "FE11"
$wsFE.Range("B13").Value = 300
$wsFE.Range("C13").Value = 11
$wsFE.Range("D13").Value = 600
$wsFE.Range("E13").Value = 40
$wsFE.Range("F13").Value = 10881
$wsFE.Range("G13").Value = 0.25
$wsFE.Range("H13").Value = $true
$wsFE.Range("I13").Value = $true
$wsFE.Range("J13").Value = $true
$wsFE.Range("K13").Value = "Todo TRUE"
$wsFE.Range("L13").Value = "Todo TRUE"

$wsFE.Range("A14").Value = "FE12"
$wsFE.Range("B14").Value = 300
$wsFE.Range("C14").Value = 11
$wsFE.Range("D14").Value = 800
$wsFE.Range("E14").Value = 40
$wsFE.Range("F14").Value = 10881
$wsFE.Range("G14").Value = 0.25
$wsFE.Range("H14").Value = $true
$wsFE.Range("I14").Value = $true
$wsFE.Range("J14").Value = $true
$wsFE.Range("K14").Value = "Todo TRUE"
$wsFE.Range("L14").Value = "Todo TRUE"

# ---------------------------------------------------------------------
# TS sheet: row 7 TS01 -> TS05
# ---------------------------------------------------------------------
$wsTS.Range("A7").Value = "TS05"

# ---------------------------------------------------------------------
# Selections per sheet, and make Sheet1 the active tab.
# ---------------------------------------------------------------------
$wsTS.Range("A8").Select() | Out-Null
$wsFE.Range("J26").Select() | Out-Null
$ws1.Range("E19").Select() | Out-Null
$ws1.Activate() | Out-Null
